$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, "Génique", "Yoann", "09/01/1996", "genique.yoann@outlook.com"),
    @(2, "Mairot", "Jean-christophe", "15/05/1999", "mairot.jean-christophe@gmail.com"),
    @(3, "Cherief", "Saufiane", "25/08/1992", "saufiane.cherief@gmail.com"),
    @(4, "Rameau", "Célia", "13/04/2000", "celia.rameau@gmail.com"),
    @(5, "Ligourel", "Teedji", "15/02/1997", "ligourel.teedji@gmail.com")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        $value = $row[$c]
        # Dates such as "09/01/1996" are ambiguous (day <= 12) and Excel's
        # input parser would silently convert them into a date serial
        # number. Force text interpretation via a temporary Text number
        # format, then strip the format back off so no residual style is
        # left on the cell (matches plain, unstyled shared-string cells).
        if ($c -eq 3) {
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.ClearFormats()
        } else {
            $cell.Value = $value
        }
    }
}
